$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Package" column (K) -------------------------------------------------
$ws.Cells.Item(1, 11).Value = "Package"

$ws.Cells.Item(4, 11).Value = "SOT-223"
$ws.Cells.Item(6, 11).Value = "8-PDIP"
$ws.Cells.Item(5, 11).Value = "3386P"
$ws.Cells.Item(3, 11).Value = "8-SOIC (0.154"", 3.90mm Width)"

# K3 picks up the same look as the existing "description" style used in B3:B5
# (Arial 9pt, black) by copying the format across.
$ws.Range("B3").Copy()
$ws.Range("K3").PasteSpecial(-4122)

# K5 uses that same font, plus wrapped / vertically centred text.
$ws.Range("B3").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("K5").WrapText = $true
$ws.Range("K5").VerticalAlignment = -4108
$ws.Application.CutCopyMode = $false

# --- Highlight the whole "link" column (G) with the built-in "Good" style ----
$ws.Range("G2").Style = "Good"
$ws.Range("G3").Style = "Good"
$ws.Range("G4").Style = "Good"
$ws.Range("G5").Style = "Good"
$ws.Range("G6").Style = "Good"
$ws.Range("G7").Style = "Good"
$ws.Range("G8").Style = "Good"
$ws.Range("G9").Style = "Good"
$ws.Range("G10").Style = "Good"
$ws.Range("G11").Style = "Good"
$ws.Range("G12").Style = "Good"
$ws.Range("G13").Style = "Good"
$ws.Range("G14").Style = "Good"
$ws.Range("G15").Style = "Good"

# G8's link text did not have a real hyperlink before - add one now, keeping
# the new green "Good" shading.
$ws.Hyperlinks.Add($ws.Range("G8"), "https://www.digikey.ca/product-detail/en/te-connectivity-amp-connectors/2-644803-2/A30924-ND/698439")

# --- New rows: female connectors ---------------------------------------------
$ws.Cells.Item(16, 1).Value = "Connector"
$ws.Cells.Item(16, 3).Value = "2 pin female connector"
$ws.Cells.Item(16, 7).Value = "https://www.digikey.ca/product-detail/en/te-connectivity-amp-connectors/3-640441-2/A30978-ND/698221"
$ws.Cells.Item(16, 8).Value = 0.24
$ws.Cells.Item(16, 9).Value = 3
$ws.Cells.Item(16, 10).Value = 3
$ws.Range("G16").Style = "Good"

$ws.Cells.Item(17, 1).Value = "Connector"
$ws.Cells.Item(17, 3).Value = "4 pin female connector"
$ws.Cells.Item(17, 7).Value = "https://www.digikey.ca/product-detail/en/te-connectivity-amp-connectors/3-640441-4/A30980-ND/698223"
$ws.Cells.Item(17, 8).Value = 0.26
$ws.Cells.Item(17, 9).Value = 3
$ws.Cells.Item(17, 10).Value = 3
$ws.Range("G17").Style = "Good"

# --- Final selection, matching the author's last action ---------------------
$ws.Range("G16:G17").Select()
